# Apply updated crypto price/volume data per commit "Updated cryptos list on Thu Aug  3 07:47:33 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "29.045.43"
$ws.Cells.Item(2, 5).Value = "  -1.87%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "1.833.33"
$ws.Cells.Item(3, 5).Value = "  -1.30%  "
# Row 4
$ws.Cells.Item(4, 4).Value = "'0.9989"
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
# Row 5
$ws.Cells.Item(5, 4).Value = "'239.61"
$ws.Cells.Item(5, 5).Value = "  -2.16%  "
# Row 6
$ws.Cells.Item(6, 4).Value = "'0.6708"
$ws.Cells.Item(6, 5).Value = "  -3.24%  "
# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
# Row 8
$ws.Cells.Item(8, 4).Value = "'0.07425"
$ws.Cells.Item(8, 5).Value = "  -3.49%  "
# Row 9
$ws.Cells.Item(9, 4).Value = "'0.2955"
$ws.Cells.Item(9, 5).Value = "  -3.52%  "
# Row 10
$ws.Cells.Item(10, 4).Value = "'22.73"
$ws.Cells.Item(10, 5).Value = "  -4.10%  "
# Row 11
$ws.Cells.Item(11, 4).Value = "'0.07646"
$ws.Cells.Item(11, 5).Value = "  -1.55%  "
# Row 12
$ws.Cells.Item(12, 4).Value = "1.843.91"
$ws.Cells.Item(12, 5).Value = "  -0.77%  "
# Row 13
$ws.Cells.Item(13, 4).Value = "'5.008"
$ws.Cells.Item(13, 5).Value = "  -2.63%  "
# Row 14
$ws.Cells.Item(14, 4).Value = "'0.6726"
$ws.Cells.Item(14, 5).Value = "  -2.73%  "
# Row 15
$ws.Cells.Item(15, 4).Value = "'86.10"
$ws.Cells.Item(15, 5).Value = "  -5.56%  "
# Row 16
$ws.Cells.Item(16, 4).Value = "'6.152"
$ws.Cells.Item(16, 5).Value = "  -6.29%  "
# Row 17
$ws.Cells.Item(17, 4).Value = "29.050.75"
$ws.Cells.Item(17, 5).Value = "  -1.45%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "'0.000008240"
$ws.Cells.Item(18, 5).Value = "  -0.71%  "
# Row 19
$ws.Cells.Item(19, 2).Value = "BitcoinCash"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(19, 4).Value = "'227.58"
$ws.Cells.Item(19, 5).Value = "  -5.05%  "
# Row 20
$ws.Cells.Item(20, 2).Value = "Avalanche"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(20, 4).Value = "'12.43"
$ws.Cells.Item(20, 5).Value = "  -2.46%  "
# Row 21
$ws.Cells.Item(21, 2).Value = "Dai"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(21, 4).Value = "'0.9998"
$ws.Cells.Item(21, 5).Value = "  +0.00%  "
# Row 22
$ws.Cells.Item(22, 2).Value = "Chainlink"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(22, 4).Value = "'7.316"
$ws.Cells.Item(22, 5).Value = "  -3.86%  "
# Row 23
$ws.Cells.Item(23, 2).Value = "BinanceUSD"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(23, 4).Value = "'0.9994"
$ws.Cells.Item(23, 5).Value = "  -0.08%  "
# Row 24
$ws.Cells.Item(24, 2).Value = "Monero"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(24, 4).Value = "'160.32"
$ws.Cells.Item(24, 5).Value = "  +0.62%  "
# Row 25
$ws.Cells.Item(25, 2).Value = "Stellar"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(25, 4).Value = "'0.1431"
$ws.Cells.Item(25, 5).Value = "  -4.26%  "
# Row 26
$ws.Cells.Item(26, 2).Value = "Cosmos"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(26, 4).Value = "'8.673"
$ws.Cells.Item(26, 5).Value = "  -2.68%  "
# Row 27
$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).Value = "'17.94"
$ws.Cells.Item(27, 5).Value = "  -1.76%  "
# Row 28
$ws.Cells.Item(28, 2).Value = "PancakeSwap"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(28, 4).Value = "'1.505"
$ws.Cells.Item(28, 5).Value = "  -2.07%  "
# Row 29
$ws.Cells.Item(29, 2).Value = "Filecoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(29, 4).Value = "'4.230"
$ws.Cells.Item(29, 5).Value = "  -0.48%  "
# Row 30
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).Value = "'4.112"
$ws.Cells.Item(30, 5).Value = "  -1.39%  "
# Row 31
$ws.Cells.Item(31, 2).Value = "Toncoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(31, 4).Value = "'1.195"
$ws.Cells.Item(31, 5).Value = "  -0.98%  "
# Row 32
$ws.Cells.Item(32, 2).Value = "Hedera"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(32, 4).Value = "'0.05378"
$ws.Cells.Item(32, 5).Value = "  +4.41%  "
# Row 33
$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(33, 4).Value = "'0.7511"
$ws.Cells.Item(33, 5).Value = "  -2.35%  "
# Row 34
$ws.Cells.Item(34, 2).Value = "LidoDAOToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(34, 4).Value = "'1.852"
$ws.Cells.Item(34, 5).Value = "  -2.08%  "
# Row 35
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(35, 4).Value = "'1.123"
$ws.Cells.Item(35, 5).Value = "  -2.44%  "
# Row 36
$ws.Cells.Item(36, 2).Value = "HuobiToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(36, 4).Value = "'2.683"
$ws.Cells.Item(36, 5).Value = "  -0.15%  "
# Row 37
$ws.Cells.Item(37, 2).Value = "Maker"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(37, 4).Value = "1.290.42"
$ws.Cells.Item(37, 5).Value = "  -3.44%  "
# Row 38
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.01810"
$ws.Cells.Item(38, 5).Value = "  -3.28%  "
# Row 39
$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(39, 4).Value = "'2.709"
$ws.Cells.Item(39, 5).Value = "  -0.67%  "
# Row 40
$ws.Cells.Item(40, 2).Value = "TrustWalletToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(40, 4).Value = "'0.9211"
$ws.Cells.Item(40, 5).Value = "  -5.02%  "
# Row 41
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).Value = "'6.027"
$ws.Cells.Item(41, 5).Value = "  +3.47%  "
# Row 42
$ws.Cells.Item(42, 2).Value = "Quant"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(42, 4).Value = "'104.22"
$ws.Cells.Item(42, 5).Value = "  -2.26%  "
# Row 43
$ws.Cells.Item(43, 2).Value = "PaxDollar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(43, 4).Value = "'0.9989"
$ws.Cells.Item(43, 5).Value = "  -0.09%  "
# Row 44
$ws.Cells.Item(44, 2).Value = "XinFinNetwork"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Cells.Item(44, 4).Value = "'0.08029"
$ws.Cells.Item(44, 5).Value = "  +22.33%  "
# Row 45
$ws.Cells.Item(45, 2).Value = "RocketPoolETH"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(45, 4).Value = "1.976.98"
$ws.Cells.Item(45, 5).Value = "  -1.34%  "
# Row 46
$ws.Cells.Item(46, 2).Value = "Mantle"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(46, 4).Value = "'0.5175"
$ws.Cells.Item(46, 5).Value = "  -0.88%  "
# Row 47
$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(47, 4).Value = "'63.51"
$ws.Cells.Item(47, 5).Value = "  +0.74%  "
# Row 48
$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 4).Value = "'1.752"
$ws.Cells.Item(48, 5).Value = "  -1.41%  "
# Row 49
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "'9.312"
$ws.Cells.Item(49, 5).Value = "  -4.77%  "
# Row 50
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.05920"
$ws.Cells.Item(50, 5).Value = "  -0.10%  "
# Row 51
$ws.Cells.Item(51, 2).Value = "Aptos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(51, 4).Value = "'6.848"
$ws.Cells.Item(51, 5).Value = "  -1.65%  "
